# Update gh-pages to output generated at 456a3b4
# Applies numeric "want-to-go" counter bumps and a venue change
# across the four worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 26786
$ws1.Range("F4").Value = 592
$ws1.Range("F5").Value = 257
$ws1.Range("D6").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws1.Range("F6").Value = 616
$ws1.Range("G6").Value = 40
$ws1.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202408/iVpZhT9M1724309121723.jpeg"
$ws1.Range("F8").Value = 555
$ws1.Range("F10").Value = 365
$ws1.Range("F11").Value = 246
$ws1.Range("F12").Value = 192
$ws1.Range("F14").Value = 304
$ws1.Range("F16").Value = 437
$ws1.Range("F17").Value = 62
$ws1.Range("F18").Value = 1562
$ws1.Range("F19").Value = 215
$ws1.Range("F20").Value = 52
$ws1.Range("F21").Value = 442

# ---- Sheet "演出" (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 235
$ws2.Range("F16").Value = 26

# ---- Sheet "本地生活" (Local Life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5094
$ws3.Range("F3").Value = 243

# ---- Sheet "全部类型" (All Types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5094
$ws4.Range("F4").Value = 243
$ws4.Range("F5").Value = 26786
$ws4.Range("F6").Value = 592
$ws4.Range("F8").Value = 257
$ws4.Range("F9").Value = 235
$ws4.Range("D10").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws4.Range("F10").Value = 616
$ws4.Range("G10").Value = 40
$ws4.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202408/iVpZhT9M1724309121723.jpeg"
$ws4.Range("F19").Value = 555
$ws4.Range("F22").Value = 365
$ws4.Range("F23").Value = 246
$ws4.Range("F24").Value = 192
$ws4.Range("F27").Value = 304
$ws4.Range("F31").Value = 437
$ws4.Range("F32").Value = 62
$ws4.Range("F34").Value = 1562
$ws4.Range("F35").Value = 215
$ws4.Range("F36").Value = 26
$ws4.Range("F37").Value = 52
$ws4.Range("F38").Value = 442
